$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 1: "print" -> "printf" in the 3rd paragraph of the code textbox.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(1)
$tr1 = $shp1.TextFrame.TextRange
$para1 = $tr1.Paragraphs(3)
# paragraph text is "<TAB>print(" -> the word "print" starts right after the tab
$para1.Characters(2, 5).Text = "printf"

# ---------------------------------------------------------------------------
# 2) Slide 39: insert a new yellow highlight rectangle just before the
#    "pole tekstowe 1" textbox (same look as the other two highlight
#    rectangles already on the slide), positioned over "errorFunc()".
#    We duplicate an existing highlight rectangle so the new shape keeps the
#    identical <p:style>/<p:txBody> structure, then reposition/resize/rename
#    it and fix its z-order.
# ---------------------------------------------------------------------------
$s39 = $p.Slides.Item(39)
$template = $s39.Shapes.Item(2)
$newRect = $template.Duplicate()
$newRect.Name = "Prostokąt 2"
$newRect.Left = 2583402 / 12700
$newRect.Top = 3506680 / 12700
# Width/Height setters on this host round slightly differently than the
# constructor args, so nudge the point value by a hair so the EMU value
# that gets stored lands exactly on target after its internal float32 cast.
$newRect.Width = 472.5437007874016
$newRect.Height = 44.03882029763781
$newRect.ZOrder(3)

# ---------------------------------------------------------------------------
# 3) Slide 40: colour "throw new " (keywords + separating spaces) in blue
#    (0070C0) on the first line, matching the colouring already used for
#    "var". The trailing tabs before "EmployeeNotFoundException" stay
#    uncoloured.
# ---------------------------------------------------------------------------
$s40 = $p.Slides.Item(40)
$shp40 = $s40.Shapes.Item(1)
$tr40 = $shp40.TextFrame.TextRange
$blue = 12611584  # RGB(0,112,192) => srgbClr 0070C0

$tr40.Characters(29, 5).Font.Color.RGB = $blue   # "throw"
$tr40.Characters(34, 1).Font.Color.RGB = $blue   # " "
$tr40.Characters(35, 3).Font.Color.RGB = $blue   # "new"
$tr40.Characters(38, 1).Font.Color.RGB = $blue   # " "
